# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest generated counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { row number -> new F-column value }
$updates = @{
    "展览" = @{
        2  = 61
        4  = 2007
        5  = 319
        8  = 2049
        9  = 10391
        12 = 268
        14 = 397
        15 = 7246
        17 = 687
        18 = 144
        20 = 270
    }
    "全部类型" = @{
        2  = 61
        4  = 2007
        5  = 319
        9  = 2049
        12 = 10391
        15 = 268
        17 = 397
        18 = 7246
        20 = 687
        21 = 144
        23 = 270
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
